# Added ability to strip out * from import headers
# This marks the "required" import columns (Name, Category, Fund) with a
# trailing " *" so the importer can recognise & later strip the marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates: flag the mandatory fields with a trailing asterisk.
$ws.Range("A1").Value = "Name *"
$ws.Range("C1").Value = "Category *"
$ws.Range("E1").Value = "Fund *"

# Bring the "Fund" column's styling in line with the rest of the header/
# data columns (drop the extra applyFill formatting it previously had).
$ws.Range("E1:E7").Style = $ws.Range("A1").Style

# Reflect the cell selection left behind after making the edits.
$ws.Range("C2").Select() | Out-Null
